# Commit: "Corrección PPT Presentacion 1"
# "Se elimina Listado de Empresas con reservas, por no considerarse util y ser redundante"
#
# Slide 17 (1-based, in sldIdLst order) contains a single full-slide picture
# showing the "Listado de Empresas con Reservas Vigentes" flow diagram.
# It is redundant with the other use-case slides, so it is deleted outright.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$s.Delete()
